$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 8).Value = 6.4   # H3: 6 -> 6.4
$ws.Cells.Item(3, 9).Value = 7.8   # I3: 7.6 -> 7.8
$ws.Cells.Item(3, 11).Value = 4.9   # K3: 4.8 -> 4.9
$ws.Cells.Item(3, 14).Value = 3.85   # N3: 3.9 -> 3.85
$ws.Cells.Item(3, 15).Value = 1.28   # O3: 1.3 -> 1.28
$ws.Cells.Item(3, 22).Value = 1.14   # V3: 1.15 -> 1.14
$ws.Cells.Item(3, 24).Value = 20   # X3: 1000 -> 20
$ws.Cells.Item(3, 26).Value = 60   # Z3: 1000 -> 60
$ws.Cells.Item(3, 27).Value = 230   # AA3: 1000 -> 230
$ws.Cells.Item(3, 28).Value = 8.6   # AB3: 970 -> 8.6
$ws.Cells.Item(3, 30).Value = 32   # AD3: 34 -> 32
$ws.Cells.Item(3, 31).Value = 120   # AE3: 1000 -> 120
$ws.Cells.Item(3, 33).Value = 970   # AG3: 12 -> 970
$ws.Cells.Item(3, 34).Value = 23   # AH3: 28 -> 23
$ws.Cells.Item(3, 35).Value = 110   # AI3: 1000 -> 110
$ws.Cells.Item(3, 38).Value = 38   # AL3: 1000 -> 38
$ws.Cells.Item(3, 39).Value = 160   # AM3: 1000 -> 160
$ws.Cells.Item(3, 40).Value = 8.800000000000001   # AN3: 10 -> 8.800000000000001
$ws.Cells.Item(3, 41).Value = 150   # AO3: 1000 -> 150

# Row 4
$ws.Cells.Item(4, 6).Value = 5.7   # F4: 5.6 -> 5.7
$ws.Cells.Item(4, 7).Value = 6   # G4: 5.9 -> 6
$ws.Cells.Item(4, 8).Value = 1.59   # H4: 1.61 -> 1.59
$ws.Cells.Item(4, 9).Value = 1.61   # I4: 1.63 -> 1.61
$ws.Cells.Item(4, 17).Value = 1.67   # Q4: 1.68 -> 1.67
$ws.Cells.Item(4, 20).Value = 1.76   # T4: 1.74 -> 1.76
$ws.Cells.Item(4, 22).Value = 2.62   # V4: 2.58 -> 2.62
$ws.Cells.Item(4, 24).Value = 24   # X4: 23 -> 24
$ws.Cells.Item(4, 27).Value = 15.5   # AA4: 16 -> 15.5
$ws.Cells.Item(4, 28).Value = 24   # AB4: 26 -> 24
$ws.Cells.Item(4, 36).Value = 160   # AJ4: 150 -> 160

# Row 5
$ws.Cells.Item(5, 6).Value = 2.1   # F5: 2.16 -> 2.1
$ws.Cells.Item(5, 7).Value = 2.3   # G5: 2.34 -> 2.3
$ws.Cells.Item(5, 8).Value = 3.4   # H5: 3.3 -> 3.4
$ws.Cells.Item(5, 9).Value = 3.9   # I5: 3.65 -> 3.9
$ws.Cells.Item(5, 10).Value = 3.45   # J5: 3.5 -> 3.45
$ws.Cells.Item(5, 11).Value = 4   # K5: 3.85 -> 4
$ws.Cells.Item(5, 14).Value = 3.8   # N5: 1.98 -> 3.8
$ws.Cells.Item(5, 15).Value = 1.28   # O5: 1.27 -> 1.28
$ws.Cells.Item(5, 17).Value = 1.84   # Q5: 1.82 -> 1.84
$ws.Cells.Item(5, 18).Value = 1.37   # R5: 1.18 -> 1.37
$ws.Cells.Item(5, 19).Value = 3.1   # S5: 1.83 -> 3.1
$ws.Cells.Item(5, 20).Value = 1.7   # T5: 1.01 -> 1.7
$ws.Cells.Item(5, 21).Value = 2.16   # U5: 2.06 -> 2.16
$ws.Cells.Item(5, 22).Value = 1.35   # V5: 1.38 -> 1.35
$ws.Cells.Item(5, 23).Value = 1.76   # W5: 1.75 -> 1.76
$ws.Cells.Item(5, 24).Value = 18.5   # X5: 1000 -> 18.5
$ws.Cells.Item(5, 25).Value = 15.5   # Y5: 20 -> 15.5
$ws.Cells.Item(5, 26).Value = 28   # Z5: 36 -> 28
$ws.Cells.Item(5, 27).Value = 70   # AA5: 1000 -> 70
$ws.Cells.Item(5, 28).Value = 11   # AB5: 15 -> 11
$ws.Cells.Item(5, 29).Value = 8.800000000000001   # AC5: 11.5 -> 8.800000000000001
$ws.Cells.Item(5, 30).Value = 16   # AD5: 20 -> 16
$ws.Cells.Item(5, 31).Value = 44   # AE5: 1000 -> 44
$ws.Cells.Item(5, 32).Value = 15   # AF5: 21 -> 15
$ws.Cells.Item(5, 33).Value = 11.5   # AG5: 15.5 -> 11.5
$ws.Cells.Item(5, 34).Value = 17.5   # AH5: 23 -> 17.5
$ws.Cells.Item(5, 35).Value = 50   # AI5: 1000 -> 50
$ws.Cells.Item(5, 36).Value = 28   # AJ5: 40 -> 28
$ws.Cells.Item(5, 37).Value = 23   # AK5: 32 -> 23
$ws.Cells.Item(5, 38).Value = 38   # AL5: 1000 -> 38
$ws.Cells.Item(5, 39).Value = 100   # AM5: 1000 -> 100
$ws.Cells.Item(5, 40).Value = 16.5   # AN5: 1000 -> 16.5
$ws.Cells.Item(5, 41).Value = 40   # AO5: 1000 -> 40

# Row 6
$ws.Cells.Item(6, 8).Value = 4.6   # H6: 4.5 -> 4.6
$ws.Cells.Item(6, 11).Value = 4.9   # K6: 5.1 -> 4.9
$ws.Cells.Item(6, 14).Value = 5.5   # N6: 5.4 -> 5.5
$ws.Cells.Item(6, 16).Value = 2.7   # P6: 2.68 -> 2.7
$ws.Cells.Item(6, 18).Value = 1.69   # R6: 1.68 -> 1.69
$ws.Cells.Item(6, 19).Value = 2.16   # S6: 2.1 -> 2.16
$ws.Cells.Item(6, 21).Value = 2.5   # U6: 2.46 -> 2.5
$ws.Cells.Item(6, 25).Value = 27   # Y6: 30 -> 27
$ws.Cells.Item(6, 26).Value = 48   # Z6: 46 -> 48
$ws.Cells.Item(6, 27).Value = 110   # AA6: 120 -> 110
$ws.Cells.Item(6, 29).Value = 12   # AC6: 13.5 -> 12
$ws.Cells.Item(6, 30).Value = 24   # AD6: 21 -> 24
$ws.Cells.Item(6, 34).Value = 17.5   # AH6: 19.5 -> 17.5
$ws.Cells.Item(6, 37).Value = 16.5   # AK6: 18.5 -> 16.5
$ws.Cells.Item(6, 38).Value = 26   # AL6: 29 -> 26
$ws.Cells.Item(6, 39).Value = 60   # AM6: 65 -> 60

# Row 7
$ws.Cells.Item(7, 6).Value = 8.4   # F7: 8.199999999999999 -> 8.4
$ws.Cells.Item(7, 8).Value = 1.37   # H7: 1.38 -> 1.37
$ws.Cells.Item(7, 9).Value = 1.38   # I7: 1.39 -> 1.38
$ws.Cells.Item(7, 10).Value = 6.4   # J7: 6.2 -> 6.4
$ws.Cells.Item(7, 16).Value = 3.55   # P7: 3.5 -> 3.55
$ws.Cells.Item(7, 19).Value = 1.94   # S7: 1.95 -> 1.94
$ws.Cells.Item(7, 20).Value = 1.64   # T7: 1.65 -> 1.64
$ws.Cells.Item(7, 21).Value = 2.44   # U7: 2.48 -> 2.44
$ws.Cells.Item(7, 22).Value = 3.6   # V7: 3.55 -> 3.6
$ws.Cells.Item(7, 25).Value = 16   # Y7: 15.5 -> 16
$ws.Cells.Item(7, 29).Value = 15.5   # AC7: 15 -> 15.5
$ws.Cells.Item(7, 32).Value = 90   # AF7: 95 -> 90
$ws.Cells.Item(7, 41).Value = 3.7   # AO7: 3.75 -> 3.7

# Row 8
$ws.Cells.Item(8, 6).Value = 2.86   # F8: 2.94 -> 2.86
$ws.Cells.Item(8, 7).Value = 3.15   # G8: 3.2 -> 3.15
$ws.Cells.Item(8, 8).Value = 2.54   # H8: 2.48 -> 2.54
$ws.Cells.Item(8, 9).Value = 2.76   # I8: 2.72 -> 2.76
$ws.Cells.Item(8, 10).Value = 3.2   # J8: 3.25 -> 3.2
$ws.Cells.Item(8, 17).Value = 2.08   # Q8: 2.12 -> 2.08
$ws.Cells.Item(8, 21).Value = 2.04   # U8: 1.98 -> 2.04
$ws.Cells.Item(8, 22).Value = 1.58   # V8: 1.6 -> 1.58
$ws.Cells.Item(8, 23).Value = 1.47   # W8: 1.46 -> 1.47
$ws.Cells.Item(8, 24).Value = 14   # X8: 12.5 -> 14
$ws.Cells.Item(8, 30).Value = 970   # AD8: 13 -> 970
$ws.Cells.Item(8, 39).Value = 120   # AM8: 130 -> 120
$ws.Cells.Item(8, 40).Value = 40   # AN8: 42 -> 40
$ws.Cells.Item(8, 41).Value = 34   # AO8: 30 -> 34

# Row 10
$ws.Cells.Item(10, 14).Value = 3.5   # N10: 1.01 -> 3.5
$ws.Cells.Item(10, 16).Value = 1.9   # P10: 1.86 -> 1.9
$ws.Cells.Item(10, 28).Value = 22   # AB10: 25 -> 22

# Row 11
$ws.Cells.Item(11, 6).Value = 1.58   # F11: 1.55 -> 1.58
$ws.Cells.Item(11, 7).Value = 1.71   # G11: 1.78 -> 1.71
$ws.Cells.Item(11, 8).Value = 5   # H11: 4.8 -> 5
$ws.Cells.Item(11, 9).Value = 7   # I11: 7.4 -> 7
$ws.Cells.Item(11, 10).Value = 3.7   # J11: 3.4 -> 3.7
$ws.Cells.Item(11, 12).Value = 1.01   # L11: 1.2 -> 1.01
$ws.Cells.Item(11, 14).Value = 5.6   # N11: 5.3 -> 5.6
$ws.Cells.Item(11, 15).Value = 1.17   # O11: 1.18 -> 1.17
$ws.Cells.Item(11, 16).Value = 2.62   # P11: 2.46 -> 2.62
$ws.Cells.Item(11, 17).Value = 1.49   # Q11: 1.54 -> 1.49
$ws.Cells.Item(11, 18).Value = 1.65   # R11: 1.61 -> 1.65
$ws.Cells.Item(11, 22).Value = 1.17   # V11: 1.16 -> 1.17
$ws.Cells.Item(11, 23).Value = 2.42   # W11: 2.28 -> 2.42

# Row 13
$ws.Cells.Item(13, 11).Value = 4.1   # K13: 4.2 -> 4.1
$ws.Cells.Item(13, 23).Value = 1.4   # W13: 1.39 -> 1.4

# Row 14
$ws.Cells.Item(14, 6).Value = 1.41   # F14: 1.35 -> 1.41
$ws.Cells.Item(14, 7).Value = 1.47   # G14: 1.58 -> 1.47
$ws.Cells.Item(14, 8).Value = 8.800000000000001   # H14: 6.8 -> 8.800000000000001
$ws.Cells.Item(14, 9).Value = 11   # I14: 13.5 -> 11
$ws.Cells.Item(14, 10).Value = 4.8   # J14: 4.4 -> 4.8
$ws.Cells.Item(14, 11).Value = 5.7   # K14: 7.4 -> 5.7
$ws.Cells.Item(14, 13).Value = 1.04   # M14: 1.01 -> 1.04
$ws.Cells.Item(14, 14).Value = 2.18   # N14: 2.16 -> 2.18
$ws.Cells.Item(14, 16).Value = 2.18   # P14: 2.16 -> 2.18
$ws.Cells.Item(14, 20).Value = 1.97   # T14: 1.71 -> 1.97
$ws.Cells.Item(14, 21).Value = 1.83   # U14: 1.64 -> 1.83
$ws.Cells.Item(14, 22).Value = 1.1   # V14: 1.08 -> 1.1
$ws.Cells.Item(14, 23).Value = 3.1   # W14: 2.72 -> 3.1
$ws.Cells.Item(14, 24).Value = 22   # X14: 29 -> 22
$ws.Cells.Item(14, 25).Value = 38   # Y14: 40 -> 38
$ws.Cells.Item(14, 26).Value = 90   # Z14: 100 -> 90
$ws.Cells.Item(14, 27).Value = 390   # AA14: 1000 -> 390
$ws.Cells.Item(14, 28).Value = 9.199999999999999   # AB14: 12.5 -> 9.199999999999999
$ws.Cells.Item(14, 29).Value = 12.5   # AC14: 16 -> 12.5
$ws.Cells.Item(14, 30).Value = 44   # AD14: 46 -> 44
$ws.Cells.Item(14, 31).Value = 180   # AE14: 1000 -> 180
$ws.Cells.Item(14, 32).Value = 9   # AF14: 12.5 -> 9
$ws.Cells.Item(14, 33).Value = 11   # AG14: 14.5 -> 11
$ws.Cells.Item(14, 35).Value = 140   # AI14: 1000 -> 140
$ws.Cells.Item(14, 36).Value = 12.5   # AJ14: 17 -> 12.5
$ws.Cells.Item(14, 37).Value = 16   # AK14: 21 -> 16
$ws.Cells.Item(14, 38).Value = 48   # AL14: 50 -> 48
$ws.Cells.Item(14, 39).Value = 170   # AM14: 1000 -> 170
$ws.Cells.Item(14, 40).Value = 6.4   # AN14: 8.6 -> 6.4
$ws.Cells.Item(14, 41).Value = 230   # AO14: 1000 -> 230

# Row 16
$ws.Cells.Item(16, 6).Value = 2.82   # F16: 2.8 -> 2.82
$ws.Cells.Item(16, 7).Value = 2.86   # G16: 2.82 -> 2.86
$ws.Cells.Item(16, 9).Value = 2.64   # I16: 2.62 -> 2.64
$ws.Cells.Item(16, 11).Value = 3.8   # K16: 3.85 -> 3.8
$ws.Cells.Item(16, 18).Value = 1.45   # R16: 1.44 -> 1.45
$ws.Cells.Item(16, 22).Value = 1.6   # V16: 1.61 -> 1.6
$ws.Cells.Item(16, 23).Value = 1.53   # W16: 1.54 -> 1.53
$ws.Cells.Item(16, 27).Value = 36   # AA16: 38 -> 36
$ws.Cells.Item(16, 36).Value = 44   # AJ16: 42 -> 44

# Row 17
$ws.Cells.Item(17, 8).Value = 1.92   # H17: 1.91 -> 1.92
$ws.Cells.Item(17, 15).Value = 1.36   # O17: 1.37 -> 1.36
$ws.Cells.Item(17, 16).Value = 1.86   # P17: 1.85 -> 1.86

# Row 18
$ws.Cells.Item(18, 6).Value = 2.98   # F18: 2.96 -> 2.98
$ws.Cells.Item(18, 7).Value = 3.05   # G18: 3 -> 3.05
$ws.Cells.Item(18, 9).Value = 2.64   # I18: 2.62 -> 2.64
$ws.Cells.Item(18, 14).Value = 3.65   # N18: 3.7 -> 3.65
$ws.Cells.Item(18, 15).Value = 1.36   # O18: 1.35 -> 1.36
$ws.Cells.Item(18, 16).Value = 1.88   # P18: 1.89 -> 1.88
$ws.Cells.Item(18, 17).Value = 2.1   # Q18: 2.08 -> 2.1
$ws.Cells.Item(18, 19).Value = 3.8   # S18: 3.75 -> 3.8
$ws.Cells.Item(18, 20).Value = 1.83   # T18: 1.81 -> 1.83
$ws.Cells.Item(18, 21).Value = 2.12   # U18: 2.14 -> 2.12
$ws.Cells.Item(18, 22).Value = 1.6   # V18: 1.61 -> 1.6
$ws.Cells.Item(18, 23).Value = 1.49   # W18: 1.5 -> 1.49
$ws.Cells.Item(18, 26).Value = 16   # Z18: 16.5 -> 16
$ws.Cells.Item(18, 38).Value = 48   # AL18: 46 -> 48
$ws.Cells.Item(18, 39).Value = 100   # AM18: 95 -> 100
$ws.Cells.Item(18, 40).Value = 32   # AN18: 34 -> 32

# Row 19
$ws.Cells.Item(19, 6).Value = 4.6   # F19: 4.5 -> 4.6
$ws.Cells.Item(19, 14).Value = 5.5   # N19: 5.4 -> 5.5
$ws.Cells.Item(19, 16).Value = 2.5   # P19: 2.48 -> 2.5
$ws.Cells.Item(19, 18).Value = 1.59   # R19: 1.58 -> 1.59
$ws.Cells.Item(19, 21).Value = 2.44   # U19: 2.48 -> 2.44
$ws.Cells.Item(19, 22).Value = 2.22   # V19: 2.2 -> 2.22
$ws.Cells.Item(19, 30).Value = 9.800000000000001   # AD19: 10 -> 9.800000000000001

# Row 20
$ws.Cells.Item(20, 6).Value = 1.58   # F20: 1.56 -> 1.58
$ws.Cells.Item(20, 7).Value = 1.6   # G20: 1.58 -> 1.6
$ws.Cells.Item(20, 8).Value = 5.9   # H20: 6 -> 5.9
$ws.Cells.Item(20, 16).Value = 3   # P20: 2.98 -> 3
$ws.Cells.Item(20, 18).Value = 1.81   # R20: 1.8 -> 1.81
$ws.Cells.Item(20, 21).Value = 2.56   # U20: 2.58 -> 2.56
$ws.Cells.Item(20, 23).Value = 2.68   # W20: 2.72 -> 2.68
$ws.Cells.Item(20, 28).Value = 14.5   # AB20: 14 -> 14.5
$ws.Cells.Item(20, 29).Value = 12   # AC20: 11.5 -> 12
$ws.Cells.Item(20, 31).Value = 60   # AE20: 65 -> 60
$ws.Cells.Item(20, 36).Value = 16.5   # AJ20: 16 -> 16.5
$ws.Cells.Item(20, 40).Value = 5.2   # AN20: 5.3 -> 5.2

# Row 21
$ws.Cells.Item(21, 6).Value = 1.31   # F21: 1.3 -> 1.31
$ws.Cells.Item(21, 10).Value = 6.6   # J21: 6.8 -> 6.6
$ws.Cells.Item(21, 19).Value = 1.98   # S21: 1.99 -> 1.98
$ws.Cells.Item(21, 21).Value = 2.28   # U21: 2.24 -> 2.28
$ws.Cells.Item(21, 23).Value = 4.1   # W21: 4.2 -> 4.1
$ws.Cells.Item(21, 35).Value = 95   # AI21: 990 -> 95
$ws.Cells.Item(21, 36).Value = 12   # AJ21: 11.5 -> 12
$ws.Cells.Item(21, 39).Value = 95   # AM21: 1000 -> 95
$ws.Cells.Item(21, 40).Value = 3.6   # AN21: 3.65 -> 3.6
$ws.Cells.Item(21, 41).Value = 95   # AO21: 110 -> 95

# Row 22
$ws.Cells.Item(22, 27).Value = 42   # AA22: 44 -> 42
